# Daily attendance processing - 2025-11-15 20:22:04
#
# The "Recorded By" column (G) lists the users who recorded each
# attendance session, separated by ", ". This pass normalizes the
# ordering of that list by swapping the first two recorded-by entries,
# while leaving single-entry cells and cells already reading
# "System, backup@backdoor.com" untouched.
#
# Note: this runtime's -eq/-ceq operators are case-insensitive, so exact
# (case-sensitive) comparisons use the .Equals() string method instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Values that must remain exactly as-is (no reordering applied to these).
$unchanged = @("System, backup@backdoor.com", "dnasr281@gmail.com", "System", "Recorded By")

$updated = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.Equals("")) { continue }

    $skip = $false
    foreach ($u in $unchanged) {
        if ($val.Equals($u)) { $skip = $true }
    }
    if ($skip) { continue }

    $parts = $val -split ", "
    if ($parts.Count -ge 2) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $newVal = $parts -join ", "
        $cell.Value = $newVal
        $updated++
    }
}

Write-Host "Updated recorded-by cells:" $updated
